$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; existing rows 9..23 shift down to 10..24.
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the new market-price record.
$ws.Cells.Item(9, 1).Value = 11
$ws.Cells.Item(9, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(9, 3).Value = "Bíobío"
$ws.Cells.Item(9, 4).Value = 44447
$ws.Cells.Item(9, 5).Value = 8
$ws.Cells.Item(9, 6).Value = 100112031
$ws.Cells.Item(9, 7).Value = "Poroto verde"
$ws.Cells.Item(9, 8).Value = "Magnum"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 37000
$ws.Cells.Item(9, 12).Value = 38000
$ws.Cells.Item(9, 13).Value = 37500
$ws.Cells.Item(9, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(9, 15).Value = "Perú"
$ws.Cells.Item(9, 16).Value = 1500
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"
